$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$ws.Range("B1").Value = "California"
$ws.Range("C1").Value = 45272
$ws.Range("C1").NumberFormat = "mm-dd-yy"
